$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to keep the value as text even if it looks numeric,
    # matching the original workbook where these "price" figures are stored
    # as plain strings rather than numbers. NumberFormat is reset back to
    # the default afterwards so the cell style does not change.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.083.27"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.173.90"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "608.23"
$ws.Range("E5").Value = "  +0.64%  "

# Row 6 - Solana
Set-TextValue "D6" "154.59"
$ws.Range("E6").Value = "  +0.26%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.171.32"
$ws.Range("E8").Value = "  -1.36%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +2.06%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.158"
$ws.Range("E10").Value = "  -1.37%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.67"
$ws.Range("E11").Value = "  -7.54%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.518"
$ws.Range("E12").Value = "  +1.32%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.47%  "

# Row 14 - Avalanche
Set-TextValue "D14" "38.41"
$ws.Range("E14").Value = "  -2.84%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.692.37"
$ws.Range("E15").Value = "  -1.37%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "66.128.34"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -1.45%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.173.42"
$ws.Range("E18").Value = "  -1.43%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +1.05%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "511.00"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -0.40%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -1.12%  "

# Row 23 - Uniswap
Set-TextValue "D23" "8.02"
$ws.Range("E23").Value = "  -1.01%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "14.87"
$ws.Range("E24").Value = "  -3.82%  "

# Row 25 - Litecoin
Set-TextValue "D25" "84.64"
$ws.Range("E25").Value = "  -0.62%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.17%  "

# Row 27 - PancakeSwap
Set-TextValue "D27" "3.01"
$ws.Range("E27").Value = "  -0.39%  "

# Row 28 - RenderToken
Set-TextValue "D28" "9.15"
$ws.Range("E28").Value = "  -1.27%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  +4.13%  "

# Row 30 - was Stacks, now NEARProtocol
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D30" "7.18"
$ws.Range("E30").Value = "  +5.09%  "

# Row 31 - was NEARProtocol, now Stacks
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D31" "2.99"
$ws.Range("E31").Value = "  +4.53%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "27.99"
$ws.Range("E32").Value = "  -0.86%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.13%  "

# Row 34 - Mantle
$ws.Range("E34").Value = "  -1.65%  "

# Row 35 - Filecoin
$ws.Range("E35").Value = "  -1.35%  "

# Row 36 - Bittensor
Set-TextValue "D36" "502.70"
$ws.Range("E36").Value = "  +3.68%  "

# Row 37 - OKB
$ws.Range("E37").Value = "  -0.30%  "

# Row 38 - Hedera
Set-TextValue "D38" "0.0883"
$ws.Range("E38").Value = "  -2.57%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.0420"
$ws.Range("E39").Value = "  -0.05%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +6.23%  "

# Row 41 - Cosmos
Set-TextValue "D41" "8.78"
$ws.Range("E41").Value = "  -1.77%  "

# Row 42 - PEPE
$ws.Range("D42").Value = "0.0₃0683"
$ws.Range("E42").Value = "  +5.98%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -4.40%  "

# Row 44 - TheGraph
Set-TextValue "D44" "0.297"
$ws.Range("E44").Value = "  -1.09%  "

# Row 45 - Fetch.AI
Set-TextValue "D45" "2.43"
$ws.Range("E45").Value = "  -0.76%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.822.57"
$ws.Range("E46").Value = "  -4.24%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "28.15"
$ws.Range("E47").Value = "  -1.97%  "

# Row 48 - ThetaToken
$ws.Range("E48").Value = "  +2.26%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  -0.11%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +0.34%  "

# Row 51 - CoreDAO
Set-TextValue "D51" "2.62"
$ws.Range("E51").Value = "  +6.42%  "
